$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 55
$ws.Range("A55").Value = 131046722
$ws.Range("B55").Value = 79275
$ws.Range("E55").Value = 185
$ws.Range("F55").Value = 'Violettgrå tagellav'
$ws.Range("G55").Value = 'Bryoria nadvornikiana'
$ws.Range("H55").Value = '(Gyeln.) Brodo & D.Hawksw.'
$ws.Range("M55").ClearContents()
$ws.Range("Q55").Value = 395391
$ws.Range("R55").Value = 6804603
$ws.Range("Z55").Value = '10:52'
$ws.Range("AB55").Value = '10:52'
$ws.Range("AC55").ClearContents()
$ws.Range("AE55").Value = $false

# Row 56
$ws.Range("A56").Value = 131046925
$ws.Range("B56").Value = 79243
$ws.Range("E56").Value = 6425
$ws.Range("F56").Value = 'Garnlav'
$ws.Range("G56").Value = 'Alectoria sarmentosa'
$ws.Range("H56").Value = '(Ach.) Ach.'
$ws.Range("M56").ClearContents()
$ws.Range("Q56").Value = 395380
$ws.Range("R56").Value = 6804774
$ws.Range("Z56").Value = '11:25'
$ws.Range("AB56").Value = '11:25'
$ws.Range("AC56").ClearContents()
$ws.Range("AE56").Value = $false

# Row 72
$ws.Range("A72").Value = 131047012
$ws.Range("B72").Value = 57884
$ws.Range("E72").Value = 100109
$ws.Range("F72").Value = 'Tretåig hackspett'
$ws.Range("G72").Value = 'Picoides tridactylus'
$ws.Range("H72").Value = '(Linnaeus, 1758)'
$ws.Range("M72").Value = 'färska spår'
$ws.Range("Q72").Value = 395446
$ws.Range("R72").Value = 6804659
$ws.Range("Z72").Value = '11:35'
$ws.Range("AB72").Value = '11:35'
$ws.Range("AC72").Value = 'Troliga spår efter tretåig hackspett (barkfälkning)'
$ws.Range("AE72").Value = $true

# Row 73
$ws.Range("A73").Value = 131046930
$ws.Range("B73").Value = 79243
$ws.Range("E73").Value = 6425
$ws.Range("F73").Value = 'Garnlav'
$ws.Range("G73").Value = 'Alectoria sarmentosa'
$ws.Range("H73").Value = '(Ach.) Ach.'
$ws.Range("M73").ClearContents()
$ws.Range("Q73").Value = 395446
$ws.Range("R73").Value = 6804802
$ws.Range("Z73").Value = '11:29'
$ws.Range("AB73").Value = '11:29'
$ws.Range("AC73").ClearContents()
$ws.Range("AE73").Value = $false

# Row 74
$ws.Range("A74").Value = 131046923
$ws.Range("B74").Value = 79243
$ws.Range("E74").Value = 6425
$ws.Range("F74").Value = 'Garnlav'
$ws.Range("G74").Value = 'Alectoria sarmentosa'
$ws.Range("H74").Value = '(Ach.) Ach.'
$ws.Range("M74").ClearContents()
$ws.Range("Q74").Value = 395364
$ws.Range("R74").Value = 6804764
$ws.Range("Z74").Value = '11:22'
$ws.Range("AB74").Value = '11:22'
$ws.Range("AC74").ClearContents()
$ws.Range("AE74").Value = $false

# Row 75
$ws.Range("A75").Value = 131046916
$ws.Range("B75").Value = 79243
$ws.Range("E75").Value = 6425
$ws.Range("F75").Value = 'Garnlav'
$ws.Range("G75").Value = 'Alectoria sarmentosa'
$ws.Range("H75").Value = '(Ach.) Ach.'
$ws.Range("M75").ClearContents()
$ws.Range("Q75").Value = 395367
$ws.Range("R75").Value = 6804698
$ws.Range("Z75").Value = '11:09'
$ws.Range("AB75").Value = '11:09'
$ws.Range("AC75").ClearContents()
$ws.Range("AE75").Value = $false

# Row 76
$ws.Range("A76").Value = 131046933
$ws.Range("B76").Value = 79243
$ws.Range("E76").Value = 6425
$ws.Range("F76").Value = 'Garnlav'
$ws.Range("G76").Value = 'Alectoria sarmentosa'
$ws.Range("H76").Value = '(Ach.) Ach.'
$ws.Range("M76").ClearContents()
$ws.Range("Q76").Value = 395458
$ws.Range("R76").Value = 6804762
$ws.Range("Z76").Value = '11:32'
$ws.Range("AB76").Value = '11:32'
$ws.Range("AC76").ClearContents()
$ws.Range("AE76").Value = $false

# Row 99
$ws.Range("A99").Value = 131046928
$ws.Range("B99").Value = 79243
$ws.Range("E99").Value = 6425
$ws.Range("F99").Value = 'Garnlav'
$ws.Range("G99").Value = 'Alectoria sarmentosa'
$ws.Range("H99").Value = '(Ach.) Ach.'
$ws.Range("M99").ClearContents()
$ws.Range("Q99").Value = 395417
$ws.Range("R99").Value = 6804797
$ws.Range("Z99").Value = '11:27'
$ws.Range("AB99").Value = '11:27'
$ws.Range("AC99").ClearContents()
$ws.Range("AE99").Value = $false

# Row 100
$ws.Range("A100").Value = 131046975
$ws.Range("B100").Value = 79243
$ws.Range("E100").Value = 6425
$ws.Range("F100").Value = 'Garnlav'
$ws.Range("G100").Value = 'Alectoria sarmentosa'
$ws.Range("H100").Value = '(Ach.) Ach.'
$ws.Range("M100").ClearContents()
$ws.Range("Q100").Value = 395535
$ws.Range("R100").Value = 6804771
$ws.Range("Z100").Value = '12:07'
$ws.Range("AB100").Value = '12:07'
$ws.Range("AC100").ClearContents()
$ws.Range("AE100").Value = $false

# Row 101
$ws.Range("A101").Value = 131046985
$ws.Range("B101").Value = 79243
$ws.Range("E101").Value = 6425
$ws.Range("F101").Value = 'Garnlav'
$ws.Range("G101").Value = 'Alectoria sarmentosa'
$ws.Range("H101").Value = '(Ach.) Ach.'
$ws.Range("M101").ClearContents()
$ws.Range("Q101").Value = 395531
$ws.Range("R101").Value = 6804664
$ws.Range("Z101").Value = '12:13'
$ws.Range("AB101").Value = '12:13'
$ws.Range("AC101").ClearContents()
$ws.Range("AE101").Value = $false
